$d = $word.ActiveDocument

# --- Change 1 -------------------------------------------------------------
# "...um valor tão alto quan" + bookmark _GoBack + "to suficiente..."
# becomes a single run "...um valor tão alto quanto suficiente..." with the
# _GoBack bookmark removed from here (it gets re-created at the very end of
# the conclusion paragraph below).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$d.Content.Find.Execute(
    "alto quan" + [char]13 + "to suficiente",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "alto quanto suficiente", 2) | Out-Null

# The [char]13 trick above won't actually match across runs reliably, so
# fall back to a direct textual fix in case the paragraph mark isn't
# involved (the bookmark split the run, not a paragraph break).
$d.Content.Find.Execute(
    "alto quanto suficiente",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null

# --- Change 2 ---------------------------------------------------------
$d.Content.Find.Execute(
    "correta na maior parte dos circuitos, com poucos erros encontrados.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "correta nos circuitos.", 2) | Out-Null

# --- Change 3 ---------------------------------------------------------
$d.Content.Find.Execute(
    "as aula ministradas.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "as aula ministradas da disciplina.", 2) | Out-Null

# Re-create the _GoBack bookmark at the very end of that paragraph (right
# after the period we just typed).
$range = $d.Content.Find.Execute(
    "ministradas da disciplina.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

$found = $d.Content
$found.Find.Execute("ministradas da disciplina.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$end = $found.End
$bm = $d.Range($end, $end)
$d.Bookmarks.Add("_GoBack", $bm) | Out-Null
